# Applies the "Add squared term and update PSM analysis" edit:
#  - Sheet1 (平衡性检验): insert a new row 5 "tertiary_share_sq", shifting the
#    former rows 5-6 (ln_fdi, ln_road_area) down to rows 6-7; also refresh the
#    matched-sample columns (C/D/G/H) for every covariate row.
#  - Sheet2 (年度统计): refresh the PS model diagnostics (columns E-I) for
#    every year row.
#  - Sheet3 (匹配概况): bump the covariate count (B5) from 5 to 6.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 平衡性检验 (balance test)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert a fresh row at position 5; this pushes the existing "ln_fdi" /
# "ln_road_area" rows down to 6 / 7 and preserves everything above it.
$ws1.Rows.Item(5).Insert()

# Updated post-match statistics for the rows that keep their position
# (ln_pgdp, ln_pop_density, tertiary_share) plus the new squared term and the
# two rows that moved down (ln_fdi, ln_road_area).
$balanceRows = @(
    @{ Row = 2; A = "ln_pgdp";            B = "38.38008244765244";  C = "0.3365521852103142";  D = "99.1231071854279";   E = "11.26926349139651";  F = "5.89907124584044e-29";  G = "0.09462444644427603"; H = "0.9246191308781114";  I = "OK" },
    @{ Row = 3; A = "ln_pop_density";     B = "35.99059577857386";  C = "3.271726598574773";   D = "90.90949586190925";  E = "10.57449205333255";  F = "9.612212033869378e-26"; G = "0.9198731486877424";  H = "0.3577092649888407";  I = "OK" },
    @{ Row = 4; A = "tertiary_share";     B = "45.20348193168135";  C = "20.489872325343";     D = "54.67191585747635";  E = "13.36117248875251";  F = "9.650858568988543e-40"; G = "5.76089804702315";    H = "9.170212484141563e-09"; I = "需检查" },
    @{ Row = 5; A = "tertiary_share_sq";  B = "45.61125931829044";  C = "21.53136706643642";   D = "52.79374569295835";  E = "13.56053058748441";  F = "7.429898583901447e-41"; G = "6.053722957041203";   H = "1.581736069367004e-09"; I = "需检查" },
    @{ Row = 6; A = "ln_fdi";             B = "32.66295094300895";  C = "-2.979673132347963";  D = "109.1224860164869";  E = "9.656641617854888";  F = "8.645295030347428e-22"; G = "-0.8377598872433861"; H = "0.4022289290092775";  I = "OK" },
    @{ Row = 7; A = "ln_road_area";       B = "-22.93534586918614"; C = "-0.7723781488215303"; D = "-96.63236755518378"; E = "-6.731473306095445"; F = "1.957702669230089e-11"; G = "-0.2171605414839901"; H = "0.8280972535816009";  I = "OK" }
)

foreach ($r in $balanceRows) {
    $ws1.Cells.Item($r.Row, 1).Value = $r.A
    $ws1.Cells.Item($r.Row, 2).Value = [double]$r.B
    $ws1.Cells.Item($r.Row, 3).Value = [double]$r.C
    $ws1.Cells.Item($r.Row, 4).Value = [double]$r.D
    $ws1.Cells.Item($r.Row, 5).Value = [double]$r.E
    $ws1.Cells.Item($r.Row, 6).Value = [double]$r.F
    $ws1.Cells.Item($r.Row, 7).Value = [double]$r.G
    $ws1.Cells.Item($r.Row, 8).Value = [double]$r.H
    $ws1.Cells.Item($r.Row, 9).Value = $r.I
}

# ---------------------------------------------------------------------------
# Sheet 2: 年度统计 (annual PS-model diagnostics) - columns E-I refreshed
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$yearRows = @(
    @{ Row = 2;  E = "1.134434208527111"; F = "0.455877884661925";  G = "0.1574095222632352"; H = "0.184749521947265";  I = "0.9351544332584775" },
    @{ Row = 3;  E = "1.136118914357281"; F = "0.455882018360749";  G = "0.1457815837570197"; H = "0.1798919439769675"; I = "0.918037859851653" },
    @{ Row = 4;  E = "1.136682560272896"; F = "0.4558875315652311"; G = "0.1415897291679757"; H = "0.1731594211844923"; I = "0.9053366706458894" },
    @{ Row = 5;  E = "1.137462449876291"; F = "0.4558802205251314"; G = "0.1295884230354892"; H = "0.2377048333801086"; I = "0.8783269748364995" },
    @{ Row = 6;  E = "1.137215553554008"; F = "0.4558831007021503"; G = "0.1340592719927949"; H = "0.2138600778478064"; I = "0.8869253444704601" },
    @{ Row = 7;  E = "1.137706560586925"; F = "0.4558579720385676"; G = "0.1343193687919772"; H = "0.2191309382299363"; I = "0.8929689733717709" },
    @{ Row = 8;  E = "1.13728077563392";  F = "0.4559031713540883"; G = "0.1408291341884509"; H = "0.2025227718341068"; I = "0.9142092367328025" },
    @{ Row = 9;  E = "1.13826111825791";  F = "0.4558808322025088"; G = "0.1353228446794389"; H = "0.2124479991246679"; I = "0.8935099677788353" },
    @{ Row = 10; E = "1.137000388450454"; F = "0.4558842820522531"; G = "0.1430123778886961"; H = "0.1664568728568378"; I = "0.9179121544253221" },
    @{ Row = 11; E = "1.133573279502847"; F = "0.4558824606103568"; G = "0.160463463870248";  H = "0.1246647451358507"; I = "0.9265988014464368" },
    @{ Row = 12; E = "1.135095788030593"; F = "0.4558822478561606"; G = "0.1569188596736608"; H = "0.1443141942799877"; I = "0.9362895668690131" },
    @{ Row = 13; E = "1.13806327037512";  F = "0.4558827905558146"; G = "0.1434022503157772"; H = "0.1763477312243915"; I = "0.920188573960699" },
    @{ Row = 14; E = "1.137277719151322"; F = "0.4558810404526684"; G = "0.1448174610017762"; H = "0.1842467028154326"; I = "0.9313249555254248" },
    @{ Row = 15; E = "1.137085861124306"; F = "0.4558829401586367"; G = "0.1452840388754392"; H = "0.1961669919200982"; I = "0.9330080285341689" },
    @{ Row = 16; E = "1.137631515809276"; F = "0.4558855760120344"; G = "0.1420477598290924"; H = "0.2182465780101087"; I = "0.9285701793552507" },
    @{ Row = 17; E = "1.139005126867447"; F = "0.4558819955097421"; G = "0.1374398636605355"; H = "0.220881263207821";  I = "0.9137200334139748" },
    @{ Row = 18; E = "1.139478091203545"; F = "0.4558811198441085"; G = "0.1367541752163105"; H = "0.2173944970758712"; I = "0.9091924839029113" }
)

foreach ($r in $yearRows) {
    $ws2.Cells.Item($r.Row, 5).Value = [double]$r.E
    $ws2.Cells.Item($r.Row, 6).Value = [double]$r.F
    $ws2.Cells.Item($r.Row, 7).Value = [double]$r.G
    $ws2.Cells.Item($r.Row, 8).Value = [double]$r.H
    $ws2.Cells.Item($r.Row, 9).Value = [double]$r.I
}

# ---------------------------------------------------------------------------
# Sheet 3: 匹配概况 (matching overview) - covariate count 5 -> 6
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(5, 2).Value = 6
